$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lucene")

# New row 15: Logistic Regression / Count Vectorizer + TFIDF + ngram(1)
# (reuses the same metric values already present in row 13)
$ws.Range("A15").Value = "Logistic Regression"
$ws.Range("B15").Value = "Count Vectorizer + TFIDF + ngram(1)"
$ws.Range("C15").Value = "0.176 0.673 0.461 0.252 0.681"
$ws.Range("D15").Value = "1.000 0.831 0.911 0.911 0.947"
$ws.Range("E15").Value = "0.097 0.523 0.300 0.144 0.517"
$ws.Range("F15").Value = "0.904 0.803 0.838 0.899 0.957"

# Blank spacer row 16 (present in sheetData with styled-but-empty cells A:E)
$ws.Range("A16:E16").Font.Bold = $false

# New row 17: Logistic Regression / Count Vectorizer + TFIDF + ngram(1) (new metrics)
$ws.Range("A17").Value = "Logistic Regression"
$ws.Range("B17").Value = "Count Vectorizer + TFIDF + ngram(1)"
$ws.Range("C17").Value = "0.170 0.683 0.515 0.252 0.668"
$ws.Range("D17").Value = "1.000 0.812 0.896 0.911 0.963"
$ws.Range("E17").Value = "0.093 0.537 0.348 0.144 0.502"
$ws.Range("F17").Value = "0.904 0.801 0.846 0.899 0.956"

# Widen column B to fit the new longer configuration label
$ws.Columns.Item(2).ColumnWidth = 29.9

# Match the recorded selection after the edit
[void]$ws.Range("E19").Select()
